# A new weekly price observation is inserted for "Femacal de La Calera -
# Zanahoria" right before the existing row 234, pushing every following
# record down by one row (old row 234 becomes 235, ... old row 309
# becomes 310). The workbook's used range grows from A1:R309 to A1:R310.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 234 - Excel shifts rows 234..309 down to
# 235..310 and the sheet's dimension/used range grows automatically.
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with the new record's data.
$ws.Cells.Item(234, 1).Value  = 3
$ws.Cells.Item(234, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(234, 3).Value  = "Coquimbo"
$ws.Cells.Item(234, 4).Value  = 44627
$ws.Cells.Item(234, 5).Value  = 5
$ws.Cells.Item(234, 6).Value  = 100114013
$ws.Cells.Item(234, 7).Value  = "Zanahoria"
$ws.Cells.Item(234, 8).Value  = "Sin especificar"
$ws.Cells.Item(234, 9).Value  = "Primera"
$ws.Cells.Item(234, 10).Value = 380
$ws.Cells.Item(234, 11).Value = 7000
$ws.Cells.Item(234, 12).Value = 7500
$ws.Cells.Item(234, 13).Value = 7237
$ws.Cells.Item(234, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(234, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(234, 16).Value = 362
$ws.Cells.Item(234, 17).Value = 20
$ws.Cells.Item(234, 18).Value = "Hortaliza"
